# edit.ps1 - apply the MLK-themed rewrite to before.docx
$d = $word.ActiveDocument

# --- Paragraph 1: Title (sz 44) ---
$p1 = $d.Paragraphs(1).Range
$p1end = $p1.End - 1
$r1 = $d.Range($p1.Start, $p1end)
$r1.Text = "The Enduring Legacy of Martin Luther King Jr.: A Beacon of Hope and Inspiration"

# --- Paragraph 2: Author (sz 36) ---
$p2 = $d.Paragraphs(2).Range
$p2end = $p2.End - 1
$r2 = $d.Range($p2.Start, $p2end)
$r2.Text = "Theodore W. Anderson"

# --- Paragraph 3: Email (sz 32) ---
$p3 = $d.Paragraphs(3).Range
$p3end = $p3.End - 1
$r3 = $d.Range($p3.Start, $p3end)
$r3.Text = "twanderson@unifiednetwork.edu"

# --- Paragraph 5: Main body (sz 24), contains manual line breaks ---
$body = @"
Martin Luther King Jr., a towering figure of the Civil Rights Movement, remains an enduring symbol of hope and inspiration for people worldwide. His unwavering commitment to nonviolent resistance and his dream of a society where all individuals are treated equally have left an indelible mark on history. In this essay, we will delve into the life and legacy of Martin Luther King Jr., exploring his pivotal role in the fight for civil rights, the impact of his message of love and equality, and his lasting influence on American society.`v`vMartin Luther King Jr., born on January 15, 1929, in Atlanta, Georgia, emerged as a prominent leader of the Civil Rights Movement in the 1950s and 1960s. He advocated for racial equality and an end to segregation through nonviolent protests, civil disobedience, and powerful oratory. Inspired by Mahatma Gandhi's teachings on nonviolence, King believed that love and compassion could overcome hatred and injustice.`v`vKing's leadership was instrumental in numerous pivotal moments of the Civil Rights Movement, including the Montgomery Bus Boycott, the March on Washington, and the Selma to Montgomery marches. His eloquent speeches, such as the iconic `"I Have a Dream`" speech, resonated with audiences across the nation, galvanizing support for civil rights and inspiring millions to join the cause.`v`vIntroduction Continued:`v`vMartin Luther King Jr.'s message of love and equality transcended racial and social boundaries, appealing to individuals from all walks of life. His emphasis on nonviolent resistance and his call for a beloved community where all people could live together in harmony and respect had a profound impact on American society. King's teachings and activism contributed to the passage of landmark legislation, such as the Civil Rights Act of 1964 and the Voting Rights Act of 1965, which outlawed discrimination and expanded voting rights for African Americans.`v`vIntroduction Concluded:`v`vMartin Luther King Jr. faced numerous challenges and adversities throughout his life. He was subjected to arrests, threats, and violence, including the infamous assassination attempt in 1968 that took his life. Despite these obstacles, King remained steadfast in his commitment to nonviolence and his pursuit of racial equality. His legacy continues to inspire individuals and movements worldwide, advocating for justice, equality, and peace for all.
"@
$p5 = $d.Paragraphs(5).Range
$p5end = $p5.End - 1
$r5 = $d.Range($p5.Start, $p5end)
$r5.Text = $body

# --- Paragraph 7: Summary body ---
$summary = @"
Martin Luther King Jr.'s life and legacy serve as a testament to the power of nonviolent resistance and the enduring impact of a dream for a better world. His unwavering commitment to equality, his message of love and compassion, and his leadership during the Civil Rights Movement have left an indelible mark on history. King's dream of a society where all individuals are treated with dignity and respect remains an aspiration for people worldwide, inspiring generations to continue the fight for justice and equality.
"@
$p7 = $d.Paragraphs(7).Range
$p7end = $p7.End - 1
$r7 = $d.Range($p7.Start, $p7end)
$r7.Text = $summary

# --- Append a new empty paragraph at the very end of the document ---
$d.Content.InsertParagraphAfter()
